$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.794.46'
$ws.Range('E2').Value = '  +0.66%  '

# Row 3
$ws.Range('D3').Value = '3.805.70'
$ws.Range('E3').Value = '  +1.44%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '597.17'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.84%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '167.68'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.43%  '

# Row 7
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('E8').Value = '  +0.35%  '

# Row 9
$ws.Range('E9').Value = '  +1.96%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.28'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.33%  '

# Row 11
$ws.Range('E11').Value = '  +0.35%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000255'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.17%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '36.03'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.11%  '

# Row 14
$ws.Range('D14').Value = '4.444.29'
$ws.Range('E14').Value = '  +1.18%  '

# Row 15
$ws.Range('D15').Value = '3.827.52'
$ws.Range('E15').Value = '  +1.83%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '18.58'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +5.65%  '

# Row 17
$ws.Range('D17').Value = '67.780.60'
$ws.Range('E17').Value = '  +0.61%  '

# Row 18
$ws.Range('E18').Value = '  +3.19%  '

# Row 19
$ws.Range('E19').Value = '  +0.17%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '461.81'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.77%  '

# Row 21
$ws.Range('E21').Value = '  -4.75%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.703'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.61%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.0000155'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.42%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '83.60'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.76%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.09'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.53%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.12'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.39%  '

# Row 27
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.17%  '

# Row 28
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.01'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.25%  '

# Row 29
$ws.Range('D29').Value = '3.955.66'
$ws.Range('E29').Value = '  +1.39%  '

# Row 30
$ws.Range('E30').Value = '  +0.45%  '

# Row 31
$ws.Range('E31').Value = '  +4.90%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.29'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.04%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '29.67'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.35%  '

# Row 34
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.08%  '

# Row 35
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.10'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.02%  '

# Row 36
$ws.Range('D36').Value = '3.745.67'
$ws.Range('E36').Value = '  +1.05%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.1000'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.60%  '

# Row 38
$ws.Range('E38').Value = '  +4.14%  '

# Row 39
$ws.Range('E39').Value = '  +0.81%  '

# Row 40
$ws.Range('E40').Value = '  +1.26%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.79'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.39%  '

# Row 42
$ws.Range('E42').Value = '  +0.03%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '48.16'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.86%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '43.67'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.67%  '

# Row 46
$ws.Range('E46').Value = '  +0.89%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.32'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.01%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '148.49'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.33%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '395.19'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.79%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '27.00'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +7.97%  '

# Row 51
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.35'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +11.72%  '
